# Rename sheet "Delivery Notes" -> "Sheet1"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Sheet1"

# Remove the "Items Count" column (F) entirely, shifting "Review Status" (G) left into F.
$ws.Range("F1:F4").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftToLeft)
